$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.055.47'
$ws.Range('E2').Value = '  +6.27%  '
$ws.Range('D3').Value = '3.118.23'
$ws.Range('E3').Value = '  +3.95%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.51'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.39%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.106.07'
$ws.Range('E8').Value = '  +3.89%  '
$ws.Range('E9').Value = '  +2.44%  '
$ws.Range('E10').Value = '  +9.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.76'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +10.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.469'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.46%  '
$ws.Range('E13').Value = '  +5.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.58'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.123'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.70%  '
$ws.Range('D16').Value = '3.633.03'
$ws.Range('E16').Value = '  +3.92%  '
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('D18').Value = '63.018.86'
$ws.Range('E18').Value = '  +6.29%  '
$ws.Range('D19').Value = '3.113.82'
$ws.Range('E19').Value = '  +3.84%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '453.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.39%  '
$ws.Range('E21').Value = '  +3.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.734'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.61%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.51'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.64'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.95%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.02'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.96%  '
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('E27').Value = '  +0.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.71'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.92%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('E30').Value = '  +5.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.85'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +11.98%  '
$ws.Range('E32').Value = '  +12.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.13'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.23%  '
$ws.Range('E34').Value = '  +4.34%  '
$ws.Range('D35').Value = '0.0₃0807'
$ws.Range('E35').Value = '  +6.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.30'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.83%  '
$ws.Range('E37').Value = '  +1.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '50.86'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.78%  '
$ws.Range('E39').Value = '  +10.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.80'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.28%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '429.46'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.66%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.950.43'
$ws.Range('E42').Value = '  +6.20%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0375'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.277'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +9.62%  '
$ws.Range('E45').Value = '  +3.53%  '
$ws.Range('E46').Value = '  +7.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '126.07'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.17%  '
$ws.Range('E48').Value = '  -0.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.76'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.80'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.23%  '
